# Edit script: applies the two substantive changes captured in the target diff.
#
# 1) The table on slide 5 gets a different built-in table style GUID.
# 2) The presentation's theme colour palette is repainted from the
#    "Integral / Red Violet" palette to the standard "Office" palette
#    (this is what the underlying theme1.xml <a:clrScheme> ends up
#    holding after the edit - the font scheme / format scheme parts of
#    the theme are already identical, only the colours - and the theme's
#    display name, which the object model does not expose a setter for -
#    change).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide  = $p.Slides.Item(5)
$tblShp = $slide.Shapes.Item(2)
$table  = $tblShp.Table
$table.ApplyStyle("{1B599A20-12F6-4E4B-9720-A0F09DB864FB}")

# --- 2. Theme colour scheme ------------------------------------------
# Order exposed by ThemeColorScheme.Item(n):
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# Values are plain VBA RGB() longs (0x00BBGGRR), i.e. R + G*256 + B*65536.
$officeThemeRGB = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
